# Generate Report for handback
# Marks 88e27452-2fe3-4ab3-8a14-1c0b33cc1c5f.md and c78a8869-5e57-4cc9-8922-f4994bf01aac.md
# as handed back (in sync) for both the zh-cn and de-de locales, and reflects the
# same status roll-up on the Overview sheet.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: rows 3 & 4 (88e27452... and c78a8869...) flip from
# "Ready for handoff" to "Handed back: in sync with en-US" for both locale
# columns.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status

# ---------------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de.
#   - Status (col B) flips to "Handed back: in sync with en-US"
#   - Latest Target File (col E) + Latest Handback File (col F) get filled in
#     (mirroring the already-handed-back row 2 pattern: target = source md,
#     handback file = the same file referenced by the handoff column).
#   - Latest Handback DateTime (col G) gets a real timestamp instead of the
#     "0001-01-01 00:00:00" placeholder.
# ---------------------------------------------------------------------------
$localeSheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-02-17 03:04:10" },
    @{ Name = "de-de"; HandbackTime = "2016-02-17 03:04:27" }
)

foreach ($locale in $localeSheets) {
    $ws = $wb.Worksheets.Item($locale.Name)

    foreach ($row in 3, 4) {
        $sourceFile = $ws.Range("A$row").Value()
        $handoffFile = $ws.Range("C$row").Value()
        $handoffDisplay = $ws.Hyperlinks.Item($ws.Range("C$row")).TextToDisplay()

        $ws.Range("B$row").Value = $status

        $ws.Hyperlinks.Add($ws.Range("E$row"), "https://github.com/OpenLocalizationTest/oltest/blob/478d0f28dba78f40cdaef6b3e4868aba2465e6fa/e2e/88e27452-2fe3-4ab3-8a14-1c0b33cc1c5f.md", "", "", "88e27452-2fe3-4ab3-8a14-1c0b33cc1c5f.md") | Out-Null
        $ws.Hyperlinks.Add($ws.Range("F$row"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c139372095957b30468156874c1d3bacd689116c/ol-handoff/OpenLocalizationTestOrg/oltest.$($locale.Name)/xinjiang/ht/88e27452-2fe3-4ab3-8a14-1c0b33cc1c5f.2cecd4163ca9d70d54ebdcf141a54d63a718dd93.$($locale.Name).xlf", "", "", "88e27452-2fe3-4ab3-8a14-1c0b33cc1c5f.2cecd4163ca9d70d54ebdcf141a54d63a718dd93.$($locale.Name).xlf") | Out-Null

        $ws.Range("G$row").Value = $locale.HandbackTime
    }
}
